$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 1

# Row 3
$ws.Range("B3").Value = 0.25
$ws.Range("C3").Value = 0.5
$ws.Range("D3").Value = 0.3333333333333333
$ws.Range("E3").Value = 2

# Row 5 (only support changes)
$ws.Range("E5").Value = 1

# Row 6
$ws.Range("B6").Value = 0
$ws.Range("C6").Value = 0
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 1

# Row 7
$ws.Range("B7").Value = 0
$ws.Range("C7").Value = 0
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 1

# Row 8
$ws.Range("B8").Value = 0.8333333333333334
$ws.Range("C8").Value = 1
$ws.Range("D8").Value = 0.9090909090909091
$ws.Range("E8").Value = 5

# Row 9 (only support changes)
$ws.Range("E9").Value = 1

# Row 10 (only support changes)
$ws.Range("E10").Value = 1

# Row 11
$ws.Range("B11").Value = 1
$ws.Range("D11").Value = 1

# Row 12 (only support changes)
$ws.Range("E12").Value = 1

# Row 13
$ws.Range("B13").Value = 0
$ws.Range("C13").Value = 0
$ws.Range("D13").Value = 0

# Row 14
$ws.Range("B14").Value = 0.6666666666666666
$ws.Range("D14").Value = 0.8
$ws.Range("E14").Value = 2

# Row 15
$ws.Range("B15").Value = 1
$ws.Range("C15").Value = 1
$ws.Range("D15").Value = 1
$ws.Range("E15").Value = 1

# Row 16
$ws.Range("B16").Value = 1
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 1

# Row 17
$ws.Range("B17").Value = 1
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = 1

# Row 18 (only support changes)
$ws.Range("E18").Value = 2

# Row 19 (only support changes)
$ws.Range("E19").Value = 1

# Row 20
$ws.Range("B20").Value = 0.6
$ws.Range("C20").Value = 1
$ws.Range("D20").Value = 0.7499999999999999
$ws.Range("E20").Value = 3

# Row 21
$ws.Range("B21").Value = 0
$ws.Range("C21").Value = 0
$ws.Range("D21").Value = 0
$ws.Range("E21").Value = 1

# Row 22
$ws.Range("B22").Value = 0
$ws.Range("C22").Value = 0
$ws.Range("D22").Value = 0
$ws.Range("E22").Value = 2

# Row 23 (accuracy row)
$ws.Range("B23").Value = 0.6774193548387096
$ws.Range("C23").Value = 0.6774193548387096
$ws.Range("D23").Value = 0.6774193548387096
$ws.Range("E23").Value = 0.6774193548387096

# Row 24 (macro avg)
$ws.Range("B24").Value = 0.5166666666666666
$ws.Range("C24").Value = 0.5952380952380952
$ws.Range("D24").Value = 0.5456709956709956
$ws.Range("E24").Value = 31

# Row 25 (weighted avg)
$ws.Range("B25").Value = 0.5580645161290323
$ws.Range("C25").Value = 0.6774193548387096
$ws.Range("D25").Value = 0.6041544477028348
$ws.Range("E25").Value = 31
